# Week 13 logging: append this week's per-play/per-unit logs and roll the
# season-to-date totals forward on the summary sheets (YDS, ST, OFF, DEF,
# TURNS, PEN).

$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 13 play-by-play logs ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "5 9 3 6 11 12 5 11 3 6 11 4 1 5 4 -2 3 2 2 0 0 -5 -1 5 3 3 5 2 0 0 0 2 4 14 3 2 3 1 4 5 6 2 4 39 12 9 5 3 1 2 11 4 10 3 3 16 2 11 1 4 2 3 -1 5 2 2 4 1 4 8 4 0 8 2 2 0 14 3 3 5 2 7 40 4 13 4 3 -2 11 4 4 14 0 0 11 2 0 6 4 3 4 -4 3 6 2 2 -2 5 0 -1 2 6 1 1 3 5 5 -2 0 2 6 8 9 6 6 0 4 5 3 4 4 14 13 3 10 -1 10 -1 5 22 1 8 4 3 4 2 2 2 12 1 2 12 3 6 8 3 3 11 1 0 0 2 6 2 3 11 4 5 3 0 0 13 6 6 6 6 20 1 3 3 10 -4 2 11 3 3 3 1 0 8 3 1 4 9 5 3 12 0 0 2 5 1 18 3 1 10 2 1 1 8 1 4 8 5 5 1 15 13 4 2 1 2 3 6 2 3 2 4 5 2 2 11 3 27 6 4 9 6 4 4 9 1 13 4 3 3 17 4 5 6 2 2 3 -1 2 4 3 3 1 15 2 -1 3 5 -4 1 2 2 3 3 4 1 3 9 47 13 2 2 10 3 0 6 -1 2 1 0 -1 2 0 2 0 1 2 2 11 2 -1 17 3 3 -2 25 10 0 8 2 0 4 11 1 2 2 3 4 -2 3 4 2 4 4 1 8 1 11 -1 8 1 9 0 2 -2 9 13 1 3 9 6 7 6 5 0 3 9 6 0 6 4 4 2 0 6 2 2 1 6 4 -1 12 7 11 5 1 6 -1 5 4 -1 2 -2 4 -2 4 5 4 18 -1 6 17 8 10 3 21 6 6 4 4 4 10 9 5 11 3 6 2 4 11 2 11 12 6 1 3 5 5 20 1 3 0 7 2 3 0 8 0 2 4 1 4 2 0 5 58 7 4 1 3 1 6 -1 9 0 5 10 12 2 9 2 10 8 9 8 5 2 1 4 3 4 8 6 7 0 5 3 1 2 4 24 1 4 9 3 12 1 2 2 5 -1 4 14 4 1 3 2 4 4 3 4 8 8 1 4 3 7 4 5 6 3 2 3 -4 -1 15 1 -2 3 0 4 11 7 -1 5 1 6 2 8 2 0 -2 3 5 5 1 2 0 66 3 5 12 4 5 -1 8 8 10 3 8 14 0 18 3 1 1 1 2 -1 1 7 10 -3 2 7 2 2 4 4 4 4 4 6 1 6 0 5 6 3 9 7 3 1 8 8 1 4 4 4 3 2 12 7 9 7 3 -1 10 2 10 4 4 2 1 3 1 -1 10 0 2 9 -1 7 1 1 -3 13 -1 1 2 -1 4 0 2 1 5 2"
$ws.Range("B3").Value = "5 4 4 2 7 6 18 4 10 12 6 15 0 6 28 8 2 22 12 11 9 45 19 17 12 7 8 14 -3 9 7 23 6 4 8 13 9 7 10 6 12 6 2 24 8 14 11 9 24 2 16 3 11 9 6 6 6 7 3 10 7 11 16 16 8 5 6 12 9 5 4 8 6 6 16 9 18 22 7 4 8 7 14 11 14 20 9 22 7 28 9 7 13 13 5 4 7 33 2 18 26 9 10 6 8 15 3 7 13 3 4 4 13 25 2 13 19 -3 11 15 7 1 11 8 1 51 12 10 5 7 6 8 12 5 35 3 9 2 3 2 10 3 7 4 16 0 8 3 51 2 0 9 5 24 8 6 14 14 26 36 36 22 5 4 1 13 9 12 6 3 28 4 73 -1 7 4 0 5 17 6 12 6 9 5 7 27 16 4 21 -1 5 9 6 10 13 12 20 12 11 6 8 1 7 5 18 12 2 7 13 21 17 7 10 9 7 8 17 15 9 -4 18 13 19 15 8 6 17 10 -1 3 3 11 -2 8 13 7 20 46 18 6 6 9 2 6 20 4 13 10 17 7 24 8 28 6 7 34 1 4 7 6 6 13 6 3 4 5 6 6 5 6 19 4 7 4 10 9 22 21 9 9 4 11 4 22 4 9 6 2 5 5 -5 2 10 18 5 1 13 4 11 22 26 27 12 12 12 6 -2 8 8 8 -2 12 10 2 3 12 11 21 8 6 17 22 2 28 10 8 9 4 13 9 0 18 5 2 4 5 14 10 7 20 3 3 6 6 8 1 14 3 19 8 34 9 3 15 6 10 18 17 9 8 5 6 21 9 2 15 3 12 9 8 6 12 21 6 13 11 10 3 13 33 15 8 3 -1 27 22 7 6 5 -3 7 15 41 7 9 6 9 5 18 0 13 9 22 18 4 4 24 7 13 25 -2 -8 6 5 19 0 6 6 13 6 0 8 14 11 20 7 9 12 13 5 12 11 6 12 2 12 9 13 6 9 8 17 7 12 10 7 7 52 7 4 27 5 20 -2 7 11 5 18 5 -2 25 7 22 13 7 13 7 18 58 9 4 8 4 5 12 19 7 8 24 1 19 4 6 8 9 11 20 10 28 11 29 11 11 21 7 14 10 7 10 20 12 6 1 7 5 17 9 10 13 7 7 10 3 9 7 0 7 6 3 6 11 9 8 11 7 4 5 9 10 9 5 8 11 4 6 12 -1 6 7 5 5 4 5 28 9 4 16 9 0 12 24 15 9 5 6 13 18 0 23 4 14 9 9 9 6 6 18 16 6 8 15 7 3 7 11 13 22 8 9 1 11 5 9 8 11 24 -1 21 -2 8 9 14 5 14 -5 21 7 11 26 20 3 8 3 28 -1 6 11 5 16 9 -2 12 9 7 17 10 9"
$ws.Range("C2").Value = "8 3 6 12 18 3 -3 1 0 9 0 0 1 8 1 3 5 2 4 6 5 1 4 3 11 1 -7 6 1 0 6 7 14 4 3 7 4 2 2 14 4 3 6 -3 1 8 -2 9 3 2 9 2 0 9 5 6 7 2 4 2 0 -2 13 29 1 7 0 2 3 5 -5 6 1 -1 1 1 0 14 4 3 0 1 1 1 1 15 0 5 2 4 2 0 3 7 0 5 13 1 2 4 11 34 8 23 8 1 4 -1 2 5 11 10 7 8 21 5 -1 0 9 5 6 2 -1 3 11 3 0 7 0 12 14 3 4 3 2 4 11 29 2 5 3 17 3 3 1 54 1 3 0 1 -2 -1 3 5 9 2 -1 1 5 2 9 5 2 1 6 19 5 5 0 0 3 3 3 8 31 6 -2 7 9 2 7 -3 2 4 0 0 0 11 12 -1 1 6 8 -2 5 14 2 3 2 3 2 7 2 18 -2 5 8 4 10 -2 2 0 9 16 1 -2 1 3 2 4 7 2 4 5 11 4 0 2 1 7 0 7 -3 4 0 0 3 6 1 4 3 5 2 5 2 2 1 20 5 4 0 5 25 3 3 0 3 2 1 1 1 3 17 1 12 1 2 2 0 3 -1 5 3 1 -1 8 4 4 7 4 7 -1 6 9 4 13 7 6 2 13 29 27 2 1 0 2 7 21 8 1 12 4 4 7 6 2 -3 1 4 4 10 0 2 4 7 1 0 10 2 1 3 12 6 10 0 3 0 1 1 7 3 4 6 5 -1 10 4 3 6 3 1 2 2 6 5 6 2 6 6 5 7 8 22 5 4 36 1 47 3 10 6 1 2 5 18 2 3 1 20 2 0 4 -2 3 6 18 2 -4 -1 3 4 11 7 1 8 2 9 3 5 2 3 -3 5 0 7 6 3 0 4 3 18 5 4 0 7 2 2 6 19 1 17 5 2 -2 2 26 1 1 -3 2 2 4 4 3 -1 8 6 12 0 3 12 6 3 3 0 23 3 5 0 2 0 0 11 56 3 2 20 6 3 1 20 5 1 13 1 -3 -3 2 5 0 14 29 -1 6 2 4 5 0 0 9 4 -2 5 45 1 6 5 6 6 5 7 -1 2 1 1 6 3 8 2 4 4 1 1 5 3 3 0 0 5 3 3 11 4 5 5 2 4 9 -4 29 4 5 0 3 3 0 3 0 5 9 3 2 7 0 1 -6 0 3 3 16 2 2 5 0 5 8 -1 3 14 2 0 3 6 5 2 5 1 1 2 0 4 7 6 6 1 9 4 2 4 4 4 0 8 11 4 0 0 -1 3 1 1 3 -3 2 2 2 7 1 0 4 3 9 4 11 5 7 0 1 4 0 1 8 1 2 4 2 4 4 4 6 -1 0 1 5 9 5 11 1 6 9 4 2 3 3 8 7 13 9 1 4 -2 0 4 15 0 9 5 4 11 10 -1 3 7 4 5 2 16 4 3 4 1 8 1 7 5 0 1 3 4 12 9 2 5 0 0 0 2 7 5 -5 4 3 -5 2 2 3 7 6 4 2 -2 -2 7 15 12 -3 1 22 4 -1 1 7 5 4 2 3 -1 34 0 8 3 7 20 5 8 4 4 9 3 1 3 6 4 0 6 0 -1 3 2 0 1 4 4 1 2 2 2 4 0 12 3 2 9 0 1 3 11 -1 25 3 5 -1 1 0 4 4 3 11 21 -2 12 0 5 2 4 4 -1 3 2 3 7 0 7 5 -2 4 19 15 1 11 7 1 27 4 12 4 2 1 2 7 0 6 1 2 6 2 5 2 7 6 0 5 9 2 2 9 4 9 -2 3 0 3 6 8 6 -1 6 7 5 0 25 2 3 2"
$ws.Range("C3").Value = "8 21 5 18 19 6 20 0 6 17 3 14 35 8 11 0 12 12 9 3 8 13 4 5 17 9 6 7 0 7 4 8 12 5 63 13 14 9 8 12 10 10 8 5 4 16 17 4 18 18 5 5 7 20 19 3 11 8 17 11 4 8 15 4 3 3 16 8 3 0 3 15 20 0 9 15 7 16 30 23 14 8 9 14 12 13 3 9 24 0 5 16 15 -3 23 36 -2 4 14 6 3 13 4 17 11 44 9 23 6 19 9 9 22 5 9 11 19 36 4 3 16 11 18 36 1 2 3 6 19 48 5 6 16 11 6 12 3 7 14 -5 27 3 4 8 10 26 6 9 11 22 5 9 9 14 11 26 10 70 8 7 14 10 3 5 8 57 19 19 9 5 7 11 5 4 6 10 7 10 12 77 6 3 7 7 5 6 10 1 6 6 78 11 -3 22 12 3 18 11 3 11 22 4 14 31 23 7 18 6 7 7 8 1 10 22 10 14 5 4 12 31 2 -2 20 23 6 -2 5 4 1 8 3 4 2 9 9 8 5 -2 13 20 21 14 11 5 17 27 4 23 15 9 13 4 5 15 14 19 19 13 5 7 14 11 12 20 3 5 1 5 10 11 11 40 12 9 12 20 5 14 4 11 5 18 9 12 9 10 13 37 11 25 6 6 17 13 5 10 5 19 9 13 9 5 17 16 0 11 10 18 28 4 27 -5 11 2 11 6 6 34 -1 44 12 3 -1 8 3 21 12 5 30 12 13 6 12 9 12 10 12 24 2 8 7 8 13 22 2 8 6 6 11 15 1 13 7 2 15 14 9 25 16 8 6 2 13 1 10 8 7 40 17 25 7 24 12 10 7 2 52 8 7 9 3 29 4 6 4 6 13 14 8 6 6 7 4 10 16 33 15 12 9 55 14 5 10 11 9 0 4 11 14 6 7 7 5 28 9 7 8 6 10 10 -2 12 5 6 34 8 7 9 18 33 13 1 36 3 7 19 7 33 -2 15 18 0 10 21 0 44 22 25 12 1 11 8 4 8 12 7 31 12 15 9 7 25 10 12 25 10 22 22 8 14 26 20 6 7 1 14 7 10 17 4 20 6 6 6 22 25 15 3 2 8 1 5 20 7 8 5 14 0 6 32 20 -1 12 18 2 12 9 8 5 9 27 1 5 8 5 9 5 8 16 7 28 5 4 8 5 28 5 9 3 7 1 16 0 4 9 9 14 1 7 2 9 11 13 7 9 5 22 8 12 15 3 5 -2 28 1 7 19 5 14 3 4 7 5 7 12 2 18 15 9 3 6 10 6 15 -5 1 -3 10 16 3 27 12 9 7 8 13 34 7 6 9 9 23 15 7 13 4 11 6 9 7 18 13 7 0 4 7 15 10 20 19 5 12 8 14 9 7 -1 27 4 35 5 9 8 7 7 7 4 19 8 7 43 29 8 3 19 22 2 5 5 3 8 1"

# --- ST sheet: append Week 13 special-teams logs ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B4").Value = "69 57 68 65 16 60 56 68 51 61 65 62 61 63 48 63 67 61 69 50 56 65 53 64 55 65 66 70 58 64 62 63 64 62 63 66 60 64 61 61 57 63 62 63 64 64"
$ws.Range("B5").Value = "41 25 24 23 2 17 12 19 12 18 13 0 16 20 8 43 18 32 25 0 14 24 5 22 23 24 23 19 24 18 20 21 45 33 28 23 27 19 26 21 24 22 26 19 65 16"
$ws.Range("B6").Value = "11 39 24 40 16 7 20 16 22 26 23 19 18 21 16 39 7 4 18 16 16 19 18 23 30 6 28 25 22 18 2 18 27 24 16 24 21 23 27 5 20 15 13 25 22 21 28 19 20 25 23 28 29 26 19 24 0"
$ws.Range("D3").Value = "31 52 53 45 35 45 63 49 46 42 52 43 48 29 47 58 47 39 40 43 51 49 30 48 56 58 43 34 59 38 43 48 67 45 55 47 62 51 40 46 59 48 47 60 44 35 41 62 61 34 44 14 50 41 46 44 38 53 41 52 61 52 44 39 46 55 52 42 49 41 58 44 63 56 64 55 44 47 56 36 43 43 57 40 42 38 45 27 37 43 49 51 54 41 53 48 33 46 42 48 46 42 41 54 59 51 52 51 39 44 42 44 68 53 31 55 40"
$ws.Range("D4").Value = "0 0 0 0 0 0 13 15 0 0 9 0 0 0 0 -5 0 6 6 0 0 0 0 0 0 0 0 0 0 0 7 0 5 0 0 0 14 0 4 0 7 8 0 0 6 0 6 0 13 0 20 32 13 0 5 11 0 0 7 0 8 0 0 14 11 0 8 0 0 0 9 0 4 0 25 27 0 0 0 0 12 0 -1 0 0 0 0 0 0 0 0 0 20 0 0 0 0 0 0 0 0 11 0 0 0 1 0 0 0 0 0 3 0 0 0 0 0"
$ws.Range("D5").Value = "0 0 0 0 8 0 0 0 0 0 13 13 10 -1 0 4 0 0 18 0 8 91 0 0 10 0 0 0 6 13 4 0 0 0 0 1 0 0 0 0 0 0 0 13 0 0 3 0 0 0 1 0 0 0 2 0 0 0 0 0 10 0 0 0 0 0 0 0 0 0 0 0 0 0 13 0 0 4 14 0 0 0 0 0 8 6 10 12 0 0 0 18 0 11 0 0 0"
# --- OFF sheet: updated season totals ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 314
$ws.Range("D2").Value = 11
$ws.Range("F2").Value = 89
$ws.Range("G2").Value = 75
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 44
$ws.Range("N2").Value = 38
$ws.Range("O2").Value = 48
$ws.Range("P2").Value = 19
$ws.Range("B3").Value = 27
$ws.Range("C3").Value = 296
$ws.Range("E3").Value = 52
$ws.Range("F3").Value = 199
$ws.Range("G3").Value = 59
$ws.Range("H3").Value = 52
$ws.Range("I3").Value = 101
$ws.Range("J3").Value = 95
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = 306
$ws.Range("Q3").Value = 833

# --- DEF sheet: updated season totals ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 389
$ws.Range("D2").Value = 31
$ws.Range("F2").Value = 113
$ws.Range("G2").Value = 112
$ws.Range("J2").Value = 63
$ws.Range("N2").Value = 18
$ws.Range("O2").Value = 29
$ws.Range("P2").Value = 21
$ws.Range("B3").Value = 15
$ws.Range("C3").Value = 286
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 61
$ws.Range("F3").Value = 172
$ws.Range("G3").Value = 61
$ws.Range("H3").Value = 40
$ws.Range("I3").Value = 91
$ws.Range("J3").Value = 100
$ws.Range("L3").Value = 431
$ws.Range("M3").Value = 298
$ws.Range("Q3").Value = 834

# --- ST sheet: updated season totals ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 105
$ws.Range("D2").Value = 117
$ws.Range("F2").Value = 18
$ws.Range("G2").Value = 16

# --- TURNS sheet: updated season totals ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("D3").Value = 18
$ws.Range("E3").Value = 9

# --- PEN sheet: updated season totals ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B3").Value = 23